# Natmi following Dr Hou advice:
# Recompute the Fn1-Itga4 LR-pairs table so every Sending-cluster/Target-cluster
# combination (ECs, FAPs, sCs) is present (rows 2-10, 3x3 grid) and the
# expression/specificity metrics reflect the updated (3-replicate) analysis.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fn1"
$ws.Range("C2").Value = "Itga4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.906497
$ws.Range("N2").Value = 68.719491
$ws.Range("O2").Value = 0.9446038650914245
$ws.Range("P2").Value = 0.9446038650914245
$ws.Range("Q2").Value = 500.501660413694
$ws.Range("R2").Value = 4504.514943723247
$ws.Range("S2").Value = 0.04741929247156782
$ws.Range("T2").Value = 0.04741929247156782

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fn1"
$ws.Range("C3").Value = "Itga4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1329193333333333
$ws.Range("N3").Value = 0.3987579999999999
$ws.Range("O3").Value = 0.005481244732096839
$ws.Range("P3").Value = 0.005481244732096839
$ws.Range("Q3").Value = 2.904256684660888
$ws.Range("R3").Value = 26.138310161948
$ws.Range("S3").Value = 0.0002751595210066011
$ws.Range("T3").Value = 0.0002751595210066011

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fn1"
$ws.Range("C4").Value = "Itga4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.210428333333333
$ws.Range("N4").Value = 3.631285
$ws.Range("O4").Value = 0.04991489017647865
$ws.Range("P4").Value = 0.04991489017647865
$ws.Range("Q4").Value = 26.44757907091222
$ws.Range("R4").Value = 238.02821163821
$ws.Range("S4").Value = 0.002505736916221005
$ws.Range("T4").Value = 0.002505736916221005

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fn1"
$ws.Range("C5").Value = "Itga4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 22.906497
$ws.Range("N5").Value = 68.719491
$ws.Range("O5").Value = 0.9446038650914245
$ws.Range("P5").Value = 0.9446038650914245
$ws.Range("Q5").Value = 8820.20325653309
$ws.Range("R5").Value = 79381.82930879781
$ws.Range("S5").Value = 0.8356571635237158
$ws.Range("T5").Value = 0.8356571635237158

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fn1"
$ws.Range("C6").Value = "Itga4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1329193333333333
$ws.Range("N6").Value = 0.3987579999999999
$ws.Range("O6").Value = 0.005481244732096839
$ws.Range("P6").Value = 0.005481244732096839
$ws.Range("Q6").Value = 51.18091765505977
$ws.Range("R6").Value = 460.628258895538
$ws.Range("S6").Value = 0.004849060642960667
$ws.Range("T6").Value = 0.004849060642960667

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fn1"
$ws.Range("C7").Value = "Itga4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.210428333333333
$ws.Range("N7").Value = 3.631285
$ws.Range("O7").Value = 0.04991489017647865
$ws.Range("P7").Value = 0.04991489017647865
$ws.Range("Q7").Value = 466.0784199114595
$ws.Range("R7").Value = 4194.705779203136
$ws.Range("S7").Value = 0.04415791326286477
$ws.Range("T7").Value = 0.04415791326286476

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fn1"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 22.906497
$ws.Range("N8").Value = 68.719491
$ws.Range("O8").Value = 0.9446038650914245
$ws.Range("P8").Value = 0.9446038650914245
$ws.Range("Q8").Value = 649.410162162063
$ws.Range("R8").Value = 5844.691459458567
$ws.Range("S8").Value = 0.06152740909614091
$ws.Range("T8").Value = 0.06152740909614091

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fn1"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1329193333333333
$ws.Range("N9").Value = 0.3987579999999999
$ws.Range("O9").Value = 0.005481244732096839
$ws.Range("P9").Value = 0.005481244732096839
$ws.Range("Q9").Value = 3.768326768360665
$ws.Range("R9").Value = 33.91494091524599
$ws.Range("S9").Value = 0.0003570245681295712
$ws.Range("T9").Value = 0.0003570245681295712

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fn1"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.210428333333333
$ws.Range("N10").Value = 3.631285
$ws.Range("O10").Value = 0.04991489017647865
$ws.Range("P10").Value = 0.04991489017647865
$ws.Range("Q10").Value = 34.31622304517167
$ws.Range("R10").Value = 308.846007406545
$ws.Range("S10").Value = 0.00325123999739288
$ws.Range("T10").Value = 0.00325123999739288
